$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.806.19'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.642.16'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.35'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  -0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.21'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0848'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '1.871.49'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('D13').Value = '1.634.61'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.526'
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.15'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '26.820.26'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '0.0₃0735'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.77'
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.64'
$ws.Range('E21').Value = '  +5.38%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.36'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.19'
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.12'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0508'
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('E32').Value = '  +1.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('D35').Value = '1.267.20'
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0175'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.532'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('D43').Value = '1.781.48'
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('E44').Value = '  -4.77%  '
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.21'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.57'
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0966'
$ws.Range('E51').Value = '  -1.32%  '
